$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells we touch to text format before assignment,
# to avoid Excel auto-converting numeric-looking strings to numbers,
# then reset style back to Normal to avoid leaving a residual style index.
$dCells = @("D2", "D3", "D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D34", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value2 = "51.288.96"
$ws.Range("E2").Value2 = "  -0.61%  "

# Row 3
$ws.Range("D3").Value2 = "3.069.41"
$ws.Range("E3").Value2 = "  +1.00%  "

# Row 4
$ws.Range("E4").Value2 = "  -0.04%  "

# Row 5
$ws.Range("D5").Value2 = "394.20"
$ws.Range("E5").Value2 = "  +2.61%  "

# Row 6
$ws.Range("D6").Value2 = "102.06"
$ws.Range("E6").Value2 = "  -0.53%  "

# Row 7
$ws.Range("E7").Value2 = "  -1.84%  "

# Row 8
$ws.Range("E8").Value2 = "  +0.01%  "

# Row 9
$ws.Range("D9").Value2 = "0.585"
$ws.Range("E9").Value2 = "  -0.72%  "

# Row 10
$ws.Range("D10").Value2 = "37.16"
$ws.Range("E10").Value2 = "  +1.00%  "

# Row 11
$ws.Range("D11").Value2 = "0.139"
$ws.Range("E11").Value2 = "  +0.65%  "

# Row 12
$ws.Range("D12").Value2 = "0.0851"
$ws.Range("E12").Value2 = "  -1.31%  "

# Row 13
$ws.Range("D13").Value2 = "3.545.66"
$ws.Range("E13").Value2 = "  +0.82%  "

# Row 14
$ws.Range("D14").Value2 = "18.59"
$ws.Range("E14").Value2 = "  -0.65%  "

# Row 15
$ws.Range("D15").Value2 = "7.69"
$ws.Range("E15").Value2 = "  -1.03%  "

# Row 16
$ws.Range("B16").Value2 = "Polygon"
$ws.Range("C16").Value2 = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value2 = "1.03"
$ws.Range("E16").Value2 = "  +5.50%  "

# Row 17
$ws.Range("B17").Value2 = "WrappedEther"
$ws.Range("C17").Value2 = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value2 = "3.062.59"
$ws.Range("E17").Value2 = "  +0.70%  "

# Row 18
$ws.Range("D18").Value2 = "10.50"
$ws.Range("E18").Value2 = "  -0.60%  "

# Row 19
$ws.Range("D19").Value2 = "51.257.80"
$ws.Range("E19").Value2 = "  -0.74%  "

# Row 20
$ws.Range("D20").Value2 = "3.15"
$ws.Range("E20").Value2 = "  +1.56%  "

# Row 21
$ws.Range("B21").Value2 = "InternetComputer(DFINITY)"
$ws.Range("C21").Value2 = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D21").Value2 = "12.32"
$ws.Range("E21").Value2 = "  -0.57%  "

# Row 22
$ws.Range("B22").Value2 = "ShibaInu"
$ws.Range("C22").Value2 = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value2 = "0.0₃0957"
$ws.Range("E22").Value2 = "  -0.67%  "

# Row 23
$ws.Range("D23").Value2 = "69.97"
$ws.Range("E23").Value2 = "  +0.02%  "

# Row 24
$ws.Range("D24").Value2 = "264.59"
$ws.Range("E24").Value2 = "  -0.91%  "

# Row 25
$ws.Range("D25").Value2 = "3.19"
$ws.Range("E25").Value2 = "  +0.85%  "

# Row 26
$ws.Range("D26").Value2 = "7.85"
$ws.Range("E26").Value2 = "  -6.49%  "

# Row 27
$ws.Range("D27").Value2 = "26.93"
$ws.Range("E27").Value2 = "  +2.00%  "

# Row 28
$ws.Range("B28").Value2 = "RenderToken"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value2 = "7.17"
$ws.Range("E28").Value2 = "  -3.13%  "

# Row 29
$ws.Range("B29").Value2 = "Dai"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value2 = "1.00"
$ws.Range("E29").Value2 = "  +0.03%  "

# Row 30
$ws.Range("E30").Value2 = "  -5.31%  "

# Row 31
$ws.Range("E31").Value2 = "  -1.51%  "

# Row 32
$ws.Range("D32").Value2 = "10.62"
$ws.Range("E32").Value2 = "  +3.36%  "

# Row 33
$ws.Range("D33").Value2 = "0.0497"
$ws.Range("E33").Value2 = "  +11.77%  "

# Row 34
$ws.Range("D34").Value2 = "36.13"
$ws.Range("E34").Value2 = "  +6.11%  "

# Row 35
$ws.Range("E35").Value2 = "  +0.29%  "

# Row 36
$ws.Range("D36").Value2 = "50.05"
$ws.Range("E36").Value2 = "  -1.15%  "

# Row 37
$ws.Range("D37").Value2 = "0.999"
$ws.Range("E37").Value2 = "  -0.14%  "

# Row 38
$ws.Range("D38").Value2 = "3.31"
$ws.Range("E38").Value2 = "  -1.26%  "

# Row 39
$ws.Range("E39").Value2 = "  +0.95%  "

# Row 40
$ws.Range("D40").Value2 = "3.94"
$ws.Range("E40").Value2 = "  +7.53%  "

# Row 41
$ws.Range("D41").Value2 = "128.69"
$ws.Range("E41").Value2 = "  +0.33%  "

# Row 42
$ws.Range("D42").Value2 = "1.84"
$ws.Range("E42").Value2 = "  -1.23%  "

# Row 43
$ws.Range("D43").Value2 = "16.60"
$ws.Range("E43").Value2 = "  -2.15%  "

# Row 44
$ws.Range("D44").Value2 = "2.54"
$ws.Range("E44").Value2 = "  +0.80%  "

# Row 45
$ws.Range("E45").Value2 = "  -0.83%  "

# Row 46
$ws.Range("D46").Value2 = "21.64"
$ws.Range("E46").Value2 = "  -0.15%  "

# Row 47
$ws.Range("D47").Value2 = "2.50"
$ws.Range("E47").Value2 = "  +0.52%  "

# Row 48
$ws.Range("E48").Value2 = "  -1.65%  "

# Row 49
$ws.Range("D49").Value2 = "2.071.64"
$ws.Range("E49").Value2 = "  +1.88%  "

# Row 50
$ws.Range("B50").Value2 = "FlareNetwork"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/2hOSU_JYX+flarenetwork-flr"
$ws.Range("D50").Value2 = "0.0481"
$ws.Range("E50").Value2 = "  +22.66%  "

# Row 51
$ws.Range("B51").Value2 = "Mantle"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value2 = "0.903"
$ws.Range("E51").Value2 = "  +10.27%  "

# Reset style on touched D cells back to Normal (removes quote/text formatting marker)
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Done updating cryptos sheet."